$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE ON ORDERING: the workbook's shared-string table only ever appends new
# unique strings (and prunes ones that become unused) when Excel re-saves the
# file, so the order below intentionally introduces each *new* string value
# ("2 Days", the updated "...over serial." sentence, "5 Days", "3 Days") in
# the same sequence the target workbook uses, to keep the saved XML identical.

# --- Row 2: Test Component Functionality ---
# C2/D2 text content is unchanged (only the shared-string index shifts
# because other strings are removed/added elsewhere), so nothing to set here.
$ws.Range("E2").Value = "2 Days"    # Duration: 1 Day -> 2 Days
$ws.Range("G2").Value = 44213       # Finish: 1/16/2021 -> 1/17/2021

# --- Row 3: Assemble Robot Structure ---
$ws.Range("E3").Value = "2 Days"
$ws.Range("F3").Value = 44214       # Start: 1/17/2021 -> 1/18/2021
$ws.Range("G3").Value = 44215       # Finish: 1/17/2021 -> 1/19/2021

# --- Row 4: Sensor Placement and Calibration ---
$ws.Range("E4").Value = "2 Days"
$ws.Range("F4").Value = 44216       # Start: 1/20/2021
$ws.Range("G4").Value = 44217       # Finish: 1/21/2021

# --- Row 6: Link Arduino and Raspberry Pi communication ---
# (Set before row 5's new Duration value so the new unique strings are
# appended to the shared string table in the same order as the target file.)
$ws.Range("D6").Value = "Able to send text based messages to each other over serial."

# --- Row 5: Motor Placement and Calibration ---
$ws.Range("E5").Value = "5 Days"
$ws.Range("F5").Value = 44218       # Start: 1/22/2021
$ws.Range("G5").Value = 44222       # Finish: 1/26/2021

# --- Row 6 (continued) ---
$ws.Range("E6").Value = "3 Days"
# F6/G6 still carry the plain (unformatted) cell style, so first copy the
# date number format from F2 (which already has the date style applied) onto
# F6:G7, reusing the existing style entry instead of creating a new one, then
# set the actual date values.
$ws.Range("F2").Copy()
$ws.Range("F6:G7").PasteSpecial(-4122)
$ws.Range("F6").Value = 44223       # Start: 1/27/2021
$ws.Range("G6").Value = 44225       # Finish: 1/29/2021

# --- Row 7: Create basic obstacle avoidance algorithm ---
$ws.Range("E7").Value = "2 Days"
$ws.Range("F7").Value = 44226       # Start: 1/30/2021
$ws.Range("G7").Value = 44227       # Finish: 1/31/2021
$ws.Range("H7").Value = 6           # Precedence

# --- New rows 12-16, matching the existing blank "Arduino" rows (8-11) ---
$ws.Range("A11:H11").Copy()
$ws.Range("A12:H16").PasteSpecial(-4122)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Arduino"
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Arduino"
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Arduino"
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Arduino"
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Arduino"
